$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Centurion Mine, Australia, M4352, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19) # Column S = build_version
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
